# Domains template: remove the "Website Name" and "Niche" columns.
#   - "Website Name" was a free-text duplicate of the URL; drop it.
#   - "Niche" is now auto-detected instead of manually entered; drop it.
# Deleting the whole column (instead of rewriting every cell) keeps each
# remaining column's width, style and data glued together as Excel shifts
# everything left, so we don't have to re-enter headers/values by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete "Niche" (column E) first so the "Website Name" delete below
# doesn't change its column letter.
$ws.Range("E1:E1").EntireColumn.Delete()

# Delete "Website Name" (column A).
$ws.Range("A1:A1").EntireColumn.Delete()

Write-Output "Removed 'Website Name' and 'Niche' columns; sheet is now A1:I4."
